$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Subcriber Newletter form") is now finished: mark Completed and
# stamp the start/end date (19-07-2025).
$ws.Range("D10").Value = "Completed"
$ws.Range("E10").Value = "19-07-2025"
$ws.Range("F10").Value = "19-07-2025"

# New row 13: another completed task describing message/form fine tuning.
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Messages fine tune for all form"
$ws.Range("C13").Value = "Medium"
$ws.Range("D13").Value = "Completed"
$ws.Range("E13").Value = "18-07-2025"
$ws.Range("F13").Value = "18-07-2025"

# Row 16 height settles back to the normal 19.5 (was 20.25).
$ws.Rows.Item(16).RowHeight = 19.5
